# Auto-generated script to update cryptos.xlsx price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.706.79"
$ws.Range("E2").Value = "  +0.58%  "
$ws.Range("D3").Value = "'2.119.09"
$ws.Range("E3").Value = "  +0.44%  "
$ws.Range("D4").Value = "'1.017"
$ws.Range("E4").Value = "  +1.43%  "
$ws.Range("D5").Value = "'338.92"
$ws.Range("E5").Value = "  +1.24%  "
$ws.Range("D6").Value = "'1.015"
$ws.Range("E6").Value = "  +1.29%  "
$ws.Range("D7").Value = "'0.5273"
$ws.Range("E7").Value = "  +0.62%  "
$ws.Range("D8").Value = "'0.4540"
$ws.Range("E8").Value = "  -0.25%  "
$ws.Range("D9").Value = "'53.54"
$ws.Range("E9").Value = "  +0.34%  "
$ws.Range("D10").Value = "'0.09092"
$ws.Range("D11").Value = "'1.173"
$ws.Range("E11").Value = "  +0.87%  "
$ws.Range("E12").Value = "  -0.83%  "
$ws.Range("D13").Value = "'2.112.67"
$ws.Range("E13").Value = "  +0.13%  "
$ws.Range("D14").Value = "'6.826"
$ws.Range("E14").Value = "  +0.57%  "
$ws.Range("D15").Value = "'8.095"
$ws.Range("E15").Value = "  +2.98%  "
$ws.Range("D16").Value = "'97.97"
$ws.Range("E16").Value = "  +1.20%  "
$ws.Range("D17").Value = "'0.00001164"
$ws.Range("E17").Value = "  +3.23%  "
$ws.Range("D18").Value = "'1.016"
$ws.Range("E18").Value = "  +1.28%  "
$ws.Range("D19").Value = "'0.06730"
$ws.Range("E19").Value = "  +1.57%  "
$ws.Range("D20").Value = "'19.50"
$ws.Range("E20").Value = "  +0.82%  "
$ws.Range("D21").Value = "'1.014"
$ws.Range("E21").Value = "  +1.23%  "
$ws.Range("D22").Value = "'6.437"
$ws.Range("E22").Value = "  +2.00%  "
$ws.Range("D23").Value = "'30.778.72"
$ws.Range("E23").Value = "  +0.67%  "
$ws.Range("D24").Value = "'12.85"
$ws.Range("E24").Value = "  +3.84%  "
$ws.Range("D25").Value = "'2.382"
$ws.Range("E25").Value = "  +1.04%  "
$ws.Range("D26").Value = "'2.363.16"
$ws.Range("E26").Value = "  +0.26%  "
$ws.Range("D27").Value = "'22.46"
$ws.Range("E27").Value = "  +0.22%  "
$ws.Range("D28").Value = "'165.51"
$ws.Range("E28").Value = "  +1.29%  "
$ws.Range("D29").Value = "'2.538"
$ws.Range("E29").Value = "  -1.56%  "
$ws.Range("D30").Value = "'136.29"
$ws.Range("E30").Value = "  +2.42%  "
$ws.Range("D31").Value = "'1.200"
$ws.Range("E31").Value = "  -0.24%  "
$ws.Range("D32").Value = "'0.1077"
$ws.Range("E32").Value = "  +0.21%  "
$ws.Range("D33").Value = "'6.384"
$ws.Range("E33").Value = "  +3.46%  "
$ws.Range("D34").Value = "'1.630"
$ws.Range("E34").Value = "  -2.19%  "
$ws.Range("D35").Value = "'3.957"
$ws.Range("E35").Value = "  +0.13%  "
$ws.Range("E36").Value = "  -0.79%  "
$ws.Range("D37").Value = "'5.956"
$ws.Range("E37").Value = "  +7.41%  "
$ws.Range("D38").Value = "'0.02662"
$ws.Range("E38").Value = "  +3.07%  "
$ws.Range("D39").Value = "'0.06865"
$ws.Range("E39").Value = "  +0.77%  "
$ws.Range("E40").Value = "  +1.32%  "
$ws.Range("D41").Value = "'12.64"
$ws.Range("E41").Value = "  -1.12%  "
$ws.Range("D42").Value = "'0.6878"
$ws.Range("E42").Value = "  -0.76%  "
$ws.Range("D43").Value = "'1.262"
$ws.Range("E43").Value = "  +0.40%  "
$ws.Range("D44").Value = "'15.19"
$ws.Range("E44").Value = "  +8.08%  "
$ws.Range("D45").Value = "'0.6450"
$ws.Range("E45").Value = "  +0.51%  "
$ws.Range("D46").Value = "'2.309"
$ws.Range("E46").Value = "  -3.51%  "
$ws.Range("D47").Value = "'0.00000000369"
$ws.Range("E47").Value = "  +15.36%  "
$ws.Range("D48").Value = "'3.706"
$ws.Range("E48").Value = "  +1.02%  "
$ws.Range("D49").Value = "'1.257"
$ws.Range("E49").Value = "  +0.52%  "
$ws.Range("D50").Value = "'0.07331"
$ws.Range("E50").Value = "  +3.71%  "
$ws.Range("D51").Value = "'83.03"
$ws.Range("E51").Value = "  -0.54%  "
